$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Range("D:D").Insert()

# Copy number formats from the (now-shifted) column E into the new column D
# for the three data blocks on the sheet (Income Statement, Balance Sheet, Cash Flow)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New values for the newly-inserted column D (most recent period)
$d_fills = @{
  "7" = 43465
  "8" = 4316500
  "9" = 3605900
  "10" = 710600
  "12" = 157200
  "13" = 0
  "14" = 1500
  "15" = "NA"
  "17" = 4059800
  "18" = 256600
  "20" = 8100
  "21" = 836700
  "22" = 78900
  "23" = 185800
  "24" = 34000
  "25" = 0
  "26" = 151800
  "27" = 149400
  "28" = 0
  "29" = -22300
  "30" = 0
  "31" = 0
  "32" = -8100
  "33" = 127100
  "34" = 0
  "35" = 127100
  "38" = 43465
  "41" = 681600
  "42" = 0
  "43" = 724500
  "44" = 230600
  "45" = 34600
  "46" = 1671200
  "47" = "NA"
  "48" = 2650400
  "49" = 25700
  "50" = 0
  "51" = 0
  "52" = 148100
  "53" = 0
  "54" = 4495400
  "57" = 530400
  "58" = 114600
  "59" = 513400
  "60" = 1158400
  "61" = 1217700
  "62" = 263400
  "63" = 0
  "64" = 0
  "65" = 0
  "66" = 2664900
  "68" = 0
  "69" = 0
  "70" = 0
  "71" = 0
  "72" = 113200
  "73" = 0
  "74" = 0
  "75" = 0
  "76" = 1830500
  "77" = 0
  "80" = 43465
  "81" = 127100
  "83" = 572000
  "84" = 0
  "85" = 0
  "86" = 0
  "87" = 0
  "88" = 0
  "89" = 663400
  "91" = -547100
  "92" = 0
  "93" = 0
  "94" = -537400
  "96" = 0
  "97" = 0
  "98" = 0
  "99" = 0
  "100" = -40600
  "101" = -200
  "102" = 85200
}

# Revised values for column E (previously column D, restated for a handful of rows)
$e_overrides = @{
  "8" = 4207000
  "9" = 3446000
  "10" = 761100
  "17" = 3806300
  "18" = 400700
  "20" = -7700
  "21" = 975000
  "23" = 307500
  "24" = 81300
  "26" = 226200
  "27" = 222000
  "32" = 7700
  "33" = 263600
  "35" = 263600
  "81" = 263600
}

foreach ($row in $d_fills.Keys) {
  $ws.Range("D$row").Value = $d_fills[$row]
}

foreach ($row in $e_overrides.Keys) {
  $ws.Range("E$row").Value = $e_overrides[$row]
}
